$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update applicant name and source (shared strings reused elsewhere stay intact)
$ws.Range("B2").Value = "Sakura Haruno"
$ws.Range("C2").Value = "HOIT Website"

# Update exam date (was 7/6/2017, now 7/7/2017 -> serial 42923)
$ws.Range("A2").Value = "7/7/2017"

# Update contact number
$ws.Range("E2").Value = 912343134

# Update the active selection on the sheet
$ws.Range("F14").Select()
